$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin box border, centered/top aligned
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2

# Apply the exact same formatting to A2 by copying B1's format - this
# reuses B1's style entry instead of creating extra cellXfs in the
# process (setting the properties individually on a 2nd range drops one
# of them due to a quirk in the style-diffing in this engine).
$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)
